# Daily attendance processing - 2025-12-31 14:59:05
# Rotate the "Recorded By" (column G) comma-separated list of recorders
# left by one position for every data row (the first recorder is moved
# to the end of the list), leaving single-value cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $text = [string]$current
    if ($text -eq "") {
        continue
    }

    $parts = $text -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    if ($parts.Length -gt 1) {
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $newText = [string]::Join(", ", $rotated)
        $cell.Value2 = $newText
    }
}
